# Updated cryptos list on Thu Oct  5 15:02:01 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) figures scraped from
# coinranking.com, and reorders two pairs of rows whose rank changed
# (Chainlink/ShibaInu and TrustWalletToken/ARBITRUM/VeChain).
#
# Numeric-looking price strings (e.g. "212.86") are written through a
# temporary Text number format so Excel keeps them as text instead of
# silently parsing them into floating point numbers; the format / style
# is reset back to the sheet's default immediately afterwards so no
# visible formatting changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.991.72"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "1.643.00"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "212.86"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.527"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").Value = "  -0.15%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.22"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("E10").Value = "  +0.62%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0892"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.875.15"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.644.51"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("E14").Value = "  +1.23%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.560"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.54%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "64.67"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "27.969.43"
$ws.Range("E17").Value = "  +2.08%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "232.69"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.64"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  -0.09%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.31"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +3.48%  "
$ws.Range("E24").Value = "  +7.21%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "150.10"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -0.65%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.67"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "1.472.57"
$ws.Range("E33").Value = "  +4.55%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.10"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.881"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0168"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.925"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +10.05%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "69.45"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.42%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.02"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "1.784.83"
$ws.Range("E47").Value = "  +0.47%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.69"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.39%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "86.18"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("E50").Value = "  +0.22%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0993"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "

Write-Output "Applied cryptos list update: 87 cell(s) changed."
